$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4-6 (extra credential rows removed)
$ws.Range("A4:B6").EntireRow.Delete()

# Update the credential pair (row 2 keeps same shared-string slot content changed)
$ws.Range("A2").Value = "mngr353180"
$ws.Range("B2").Value = "nerynYt"

# Row 3 now uses the same new credential pair
$ws.Range("A3").Value = "mngr353180"
$ws.Range("B3").Value = "nerynYt"

# Update selection to match the new active cell state
$ws.Range("A3").Select()
